$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = -24461.9286412554
$ws.Range("B7").Value = 8041589.773721423
$ws.Range("B8").Value = 22317222.39794713
$ws.Range("B10").Value = 4096725.403303645

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("E2").Value = 120036.5755573567
$ws.Range("F2").Value = 120036.5755573567
$ws.Range("G2").Value = 120036.5755573567
$ws.Range("H2").Value = 120036.5755573567
$ws.Range("I2").Value = 120036.5755573567
$ws.Range("J2").Value = 120036.5755573567
$ws.Range("K2").Value = 120036.5755573567
$ws.Range("L2").Value = 120036.5755573567
$ws.Range("M2").Value = 120036.5755573567
$ws.Range("N2").Value = 120036.5755573567
$ws.Range("O2").Value = 120036.5755573567
$ws.Range("P2").Value = 120036.5755573567
$ws.Range("E4").Value = 66894.91200377973
$ws.Range("F4").Value = 66894.91200377973
$ws.Range("G4").Value = 66894.91200377973
$ws.Range("H4").Value = 66894.91200377973
$ws.Range("I4").Value = 66894.91200377973
$ws.Range("J4").Value = 66894.91200377973
$ws.Range("K4").Value = 66894.91200377973
$ws.Range("L4").Value = 66894.91200377973
$ws.Range("M4").Value = 66894.91200377973
$ws.Range("N4").Value = 66894.91200377973
$ws.Range("O4").Value = 66894.91200377973
$ws.Range("P4").Value = 66894.91200377973
$ws.Range("B6").Value = -48778.27397003479
$ws.Range("C6").Value = -48778.27397003479
$ws.Range("D6").Value = -48778.27397003479
$ws.Range("E6").Value = -98379.96539054526
$ws.Range("F6").Value = 34720.03460945479
$ws.Range("G6").Value = 34720.03460945479
$ws.Range("H6").Value = 34720.03460945479
$ws.Range("I6").Value = 34720.03460945479
$ws.Range("J6").Value = 34720.03460945479
$ws.Range("K6").Value = 34720.03460945479
$ws.Range("L6").Value = 34720.03460945479
$ws.Range("M6").Value = 34720.03460945479
$ws.Range("N6").Value = 34720.03460945479
$ws.Range("O6").Value = 34720.03460945479
$ws.Range("P6").Value = 34720.03460945479

$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("N11").Value = 110.5750244233121
$ws.Range("O11").Value = 117.8828208804077
$ws.Range("L12").Value = 61.18167021676314
$ws.Range("M12").Value = 51.84373129681028
$ws.Range("N12").Value = 38.66169381481656
$ws.Range("O12").Value = 57.81213424001893
$ws.Range("K14").Value = 135.370731907559
$ws.Range("L14").Value = 130.6648563030561
$ws.Range("M14").Value = 113.4004983079896
$ws.Range("N14").Value = 110.5750244233121
$ws.Range("O14").Value = 117.8828208804077
$ws.Range("P14").Value = 135.4597561231036
$ws.Range("K15").Value = 80.29914934735042
$ws.Range("L15").Value = 61.18167021676314
$ws.Range("M15").Value = 51.84373129681028
$ws.Range("N15").Value = 38.66169381481656
$ws.Range("O15").Value = 57.81213424001893
$ws.Range("P15").Value = 65.92768427608706
$ws.Range("Q15").Value = 94.49434172313325
$ws.Range("M16").Value = 92.09541281912071
$ws.Range("N16").Value = 81.96869489115805
$ws.Range("O16").Value = 96.22962838366004
$ws.Range("J17").Value = 124.5190384721106
$ws.Range("K17").Value = 135.370731907559
$ws.Range("L17").Value = 130.6648563030561
$ws.Range("M17").Value = 113.4004983079896
$ws.Range("N17").Value = 110.5750244233121
$ws.Range("O17").Value = 117.8828208804077
$ws.Range("P17").Value = 135.4597561231036
$ws.Range("Q17").Value = 150.3839754851235
$ws.Range("R17").Value = 65.71641987298243
$ws.Range("I18").Value = 10.12574714858493
$ws.Range("K18").Value = 80.29914934735042
$ws.Range("L18").Value = 61.18167021676314
$ws.Range("M18").Value = 51.84373129681028
$ws.Range("N18").Value = 38.66169381481656
$ws.Range("O18").Value = 57.81213424001893
$ws.Range("P18").Value = 65.92768427608706
$ws.Range("Q18").Value = 94.49434172313325
$ws.Range("R18").Value = 45.52166981132082
$ws.Range("K19").Value = 94.30397654773019
$ws.Range("L19").Value = 90.4687457914608
$ws.Range("M19").Value = 92.09541281912071
$ws.Range("N19").Value = 81.96869489115805
$ws.Range("O19").Value = 96.22962838366004
$ws.Range("P19").Value = 101.5955875616828
$ws.Range("Q19").Value = 65.34295837775146
$ws.Range("J20").Value = 124.5190384721106
$ws.Range("K20").Value = 135.370731907559
$ws.Range("L20").Value = 130.6648563030561
$ws.Range("M20").Value = 113.4004983079896
$ws.Range("N20").Value = 110.5750244233121
$ws.Range("O20").Value = 117.8828208804077
$ws.Range("P20").Value = 135.4597561231036
$ws.Range("Q20").Value = 150.3839754851235
$ws.Range("R20").Value = 65.71641987298243
$ws.Range("J21").Value = 93.17061249236157
$ws.Range("K21").Value = 80.29914934735042
$ws.Range("L21").Value = 61.18167021676314
$ws.Range("M21").Value = 51.84373129681028
$ws.Range("N21").Value = 38.66169381481656
$ws.Range("O21").Value = 57.81213424001893
$ws.Range("P21").Value = 65.92768427608706
$ws.Range("Q21").Value = 94.49434172313325
$ws.Range("K22").Value = 94.30397654773019
$ws.Range("L22").Value = 90.4687457914608
$ws.Range("M22").Value = 92.09541281912071
$ws.Range("N22").Value = 81.96869489115805
$ws.Range("O22").Value = 96.22962838366004
$ws.Range("P22").Value = 101.5955875616828
$ws.Range("J23").Value = 124.5190384721106
$ws.Range("K23").Value = 135.370731907559
$ws.Range("L23").Value = 130.6648563030561
$ws.Range("M23").Value = 113.4004983079896
$ws.Range("N23").Value = 110.5750244233121
$ws.Range("O23").Value = 117.8828208804077
$ws.Range("Q23").Value = 150.3839754851235
$ws.Range("L24").Value = 61.18167021676314
$ws.Range("M24").Value = 51.84373129681028
$ws.Range("N24").Value = 38.66169381481656
$ws.Range("O24").Value = 57.81213424001893
$ws.Range("P24").Value = 65.92768427608706
$ws.Range("Q24").Value = 94.49434172313325
$ws.Range("M25").Value = 92.09541281912071
$ws.Range("O25").Value = 96.22962838366004
$ws.Range("P25").Value = 101.5955875616828
$ws.Range("J26").Value = 124.5190384721106
$ws.Range("K26").Value = 135.370731907559
$ws.Range("L26").Value = 130.6648563030561
$ws.Range("M26").Value = 113.4004983079896
$ws.Range("N26").Value = 110.5750244233121
$ws.Range("O26").Value = 117.8828208804077
$ws.Range("P26").Value = 135.4597561231036
$ws.Range("Q26").Value = 150.3839754851235
$ws.Range("J27").Value = 93.17061249236157
$ws.Range("K27").Value = 80.29914934735042
$ws.Range("L27").Value = 61.18167021676314
$ws.Range("M27").Value = 51.84373129681028
$ws.Range("N27").Value = 38.66169381481656
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("Q27").Value = 94.49434172313325
$ws.Range("J28").Value = 33.63624132272333
$ws.Range("K28").Value = 94.30397654773019
$ws.Range("L28").Value = 90.4687457914608
$ws.Range("M28").Value = 92.09541281912071
$ws.Range("N28").Value = 81.96869489115805
$ws.Range("O28").Value = 96.22962838366004
$ws.Range("P28").Value = 101.5955875616828
$ws.Range("K29").Value = 135.370731907559
$ws.Range("M29").Value = 113.4004983079896
$ws.Range("N29").Value = 110.5750244233121
$ws.Range("L30").Value = 61.18167021676314
$ws.Range("M30").Value = 51.84373129681028
$ws.Range("N30").Value = 38.66169381481656
$ws.Range("O30").Value = 57.81213424001893
$ws.Range("P30").Value = 65.92768427608706
$ws.Range("Q30").Value = 94.49434172313325
$ws.Range("J32").Value = 124.5190384721106
$ws.Range("K32").Value = 135.370731907559
$ws.Range("L32").Value = 130.6648563030561
$ws.Range("M32").Value = 113.4004983079896
$ws.Range("N32").Value = 110.5750244233121
$ws.Range("O32").Value = 117.8828208804077
$ws.Range("P32").Value = 135.4597561231036
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("R32").Value = 65.71641987298243
$ws.Range("L33").Value = 61.18167021676314
$ws.Range("M33").Value = 51.84373129681028
$ws.Range("N33").Value = 38.66169381481656
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("J35").Value = 124.5190384721106
$ws.Range("K35").Value = 135.370731907559
$ws.Range("L35").Value = 130.6648563030561
$ws.Range("M35").Value = 113.4004983079896
$ws.Range("N35").Value = 110.5750244233121
$ws.Range("O35").Value = 117.8828208804077
$ws.Range("P35").Value = 135.4597561231036
$ws.Range("Q35").Value = 150.3839754851235
$ws.Range("R35").Value = 65.71641987298243
$ws.Range("M36").Value = 51.84373129681028
$ws.Range("N36").Value = 38.66169381481656
$ws.Range("L37").Value = 90.4687457914608
$ws.Range("L38").Value = 130.6648563030561
$ws.Range("M38").Value = 113.4004983079896
$ws.Range("O38").Value = 117.8828208804077
$ws.Range("M39").Value = 51.84373129681028
$ws.Range("N39").Value = 38.66169381481656
$ws.Range("K41").Value = 135.370731907559
$ws.Range("L41").Value = 130.6648563030561
$ws.Range("M41").Value = 113.4004983079896
$ws.Range("N41").Value = 110.5750244233121
$ws.Range("O41").Value = 117.8828208804077
$ws.Range("P41").Value = 135.4597561231036
$ws.Range("I42").Value = 10.12574714858493
$ws.Range("J42").Value = 93.17061249236157
$ws.Range("K42").Value = 80.29914934735042
$ws.Range("L42").Value = 61.18167021676314
$ws.Range("M42").Value = 51.84373129681028
$ws.Range("N42").Value = 38.66169381481656
$ws.Range("O42").Value = 57.81213424001893
$ws.Range("P42").Value = 65.92768427608706
$ws.Range("Q42").Value = 94.49434172313325
$ws.Range("R42").Value = 45.52166981132082
$ws.Range("M43").Value = 92.09541281912071
$ws.Range("N43").Value = 81.96869489115805
$ws.Range("O43").Value = 96.22962838366004
$ws.Range("K44").Value = 135.370731907559
$ws.Range("L44").Value = 130.6648563030561
$ws.Range("M44").Value = 113.4004983079896
$ws.Range("N44").Value = 110.5750244233121
$ws.Range("O44").Value = 117.8828208804077
$ws.Range("P44").Value = 135.4597561231036
$ws.Range("R44").Value = 65.71641987298243
$ws.Range("K45").Value = 80.29914934735042
$ws.Range("L45").Value = 61.18167021676314
$ws.Range("M45").Value = 51.84373129681028
$ws.Range("N45").Value = 38.66169381481656
$ws.Range("O45").Value = 57.81213424001893
$ws.Range("P45").Value = 65.92768427608706
$ws.Range("Q45").Value = 94.49434172313325
$ws.Range("J46").Value = 33.63624132272333
$ws.Range("K46").Value = 94.30397654773019
$ws.Range("L46").Value = 90.4687457914608
$ws.Range("M46").Value = 92.09541281912071
$ws.Range("O46").Value = 96.22962838366004
$ws.Range("P46").Value = 101.5955875616828
$ws.Range("Q46").Value = 65.34295837775146

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 108.0327934026353
$ws.Range("I18").Value = 77.12765456497084
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 78.03303713061706
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 61.14583096471014
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 108.0327934026353
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("J28").Value = 72.23757736389061
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = 108.0327934026353
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 108.0327934026353
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("I42").Value = 77.12765456497084
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 78.03303713061706
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("R44").Value = 108.0327934026353
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("J46").Value = 72.23757736389061
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 61.14583096471014

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 346980.2549709456
$ws.Range("B6").Value = 346980.2549709456
$ws.Range("B7").Value = 346980.2549709456
$ws.Range("B8").Value = 346980.2549709456
$ws.Range("B9").Value = 346980.2549709456
$ws.Range("B10").Value = 346980.2549709456
$ws.Range("B11").Value = 346980.2549709456
$ws.Range("B12").Value = 346980.2549709456
$ws.Range("B13").Value = 346980.2549709456
$ws.Range("B14").Value = 346980.2549709456
$ws.Range("B15").Value = 346980.2549709456
$ws.Range("B16").Value = 346980.2549709456

Write-Output "Applied all changes"